# Generate Report for Handback
#
# Updates the localization-status workbook after a handback:
#   - Overview/zh-cn/de-de "Status" cells move from "Ready for handoff" to
#     "Handed back: in sync with en-US".
#   - The zh-cn and de-de detail sheets gain "Latest Target File" /
#     "Latest Handback File" / "Latest Handback DateTime" data (with a
#     hyperlink on the newly-populated target-file cell), now that the
#     handback xliffs have come back.
#   - A few columns are widened so the new long file names are readable.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

$ghBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f9d2fd44db34ec939838d34bc1ded4628d6c4fcc/e2e/"
$mdA = "29a8db48-3038-46fa-a4f9-36e163ec02d9.md"
$mdB = "457d907f-d578-46a9-a135-862072b6db8b.md"

# ColumnWidth (Excel "characters" units) is stored internally with a fixed
# +5/6 bias and snapped to the nearest 1/6 step, so compensate to land as
# close as possible on the desired stored width.
function Set-ColWidth($ws, $colIndex, $targetStoredWidth) {
    $ws.Columns.Item($colIndex).ColumnWidth = $targetStoredWidth - (5.0 / 6.0)
}

# ---------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText
$wsOverview.Range("E3").Value = $statusText
$wsOverview.Range("F3").Value = $statusText

Set-ColWidth $wsOverview 5 29.9777047293527
Set-ColWidth $wsOverview 6 29.9777047293527

# ---------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = $statusText
$wsZh.Range("C3").Value = $statusText

Set-ColWidth $wsZh 3 29.9777047293527
Set-ColWidth $wsZh 9 40
Set-ColWidth $wsZh 10 40

# Row 2 (29a8db48...): Latest Target File / Latest Handback File
$wsZh.Range("J2").Value = "29a8db48-3038-46fa-a4f9-36e163ec02d9.41322cac0cab4eb5d4073b0ed8fb1c0fd8bcc28c.zh-cn.xlf"

# Row 3 (457d907f...): Latest Target File / Latest Handback File
$wsZh.Range("J3").Value = "457d907f-d578-46a9-a135-862072b6db8b.a7ed9e9cebfe4bd18a6c94ce7ec14bf903368b74.zh-cn.xlf"

# Latest Handback DateTime for both rows (zh-cn finished first)
$wsZh.Range("K2").Value = "2016-09-03 19:15:10"
$wsZh.Range("K3").Value = "2016-09-03 19:15:10"

# Recreate the hyperlinks so the new "Latest Target File" cells (I2/I3) get
# a link to the source .md file, same as column A.
$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $ghBase + $mdA, "", "", $mdA)
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $ghBase + $mdA, "", "", $mdA)
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $ghBase + $mdB, "", "", $mdB)
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $ghBase + $mdB, "", "", $mdB)

$wsZh.Range("I2").Style = "HyperLink"
$wsZh.Range("I3").Style = "HyperLink"

# ---------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = $statusText
$wsDe.Range("C3").Value = $statusText

Set-ColWidth $wsDe 3 29.9777047293527
Set-ColWidth $wsDe 9 40
Set-ColWidth $wsDe 10 40

# Row 2 (29a8db48...): Latest Target File / Latest Handback File / Latest Handback DateTime
$wsDe.Range("J2").Value = "29a8db48-3038-46fa-a4f9-36e163ec02d9.41322cac0cab4eb5d4073b0ed8fb1c0fd8bcc28c.de-de.xlf"
$wsDe.Range("K2").Value = "2016-09-03 19:15:18"

# Row 3 (457d907f...): Latest Target File / Latest Handback File / Latest Handback DateTime
$wsDe.Range("J3").Value = "457d907f-d578-46a9-a135-862072b6db8b.a7ed9e9cebfe4bd18a6c94ce7ec14bf903368b74.de-de.xlf"
$wsDe.Range("K3").Value = "2016-09-03 19:15:18"

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $ghBase + $mdA, "", "", $mdA)
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $ghBase + $mdA, "", "", $mdA)
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $ghBase + $mdB, "", "", $mdB)
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $ghBase + $mdB, "", "", $mdB)

$wsDe.Range("I2").Style = "HyperLink"
$wsDe.Range("I3").Style = "HyperLink"
